$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Plain value fixes (sign corrections for FuzzDataInserter.exe bug)
$ws.Range("D14").Value = 45752811.060000002
$ws.Range("D16").Value = 50601311.960000001
$ws.Range("D19").Value = 383100000
$ws.Range("D22").Value = 32201025

# Restore missing SUM formulas (previously hard-coded/blank values)
$ws.Range("D18").Formula = "=SUM(D12:D17)"
$ws.Range("D21").Formula = "=SUM(D18:D20)"

$excel.CalculateFullRebuild()
$wb.Save()
